$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 218, shifting existing rows 218..323 down to 219..324.
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218 with the new data record.
$ws.Cells.Item(218, 1).Value = 7
$ws.Cells.Item(218, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(218, 3).Value = "Ñuble"
$ws.Cells.Item(218, 4).Value = 44960
$ws.Cells.Item(218, 5).Value = 16
$ws.Cells.Item(218, 6).Value = 100112006
$ws.Cells.Item(218, 7).Value = "Repollo"
$ws.Cells.Item(218, 8).Value = "Crespo record"
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 200
$ws.Cells.Item(218, 11).Value = 1200
$ws.Cells.Item(218, 12).Value = 1200
$ws.Cells.Item(218, 13).Value = 1200
$ws.Cells.Item(218, 14).Value = "`$/unidad"
$ws.Cells.Item(218, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(218, 16).Value = 1200
$ws.Cells.Item(218, 17).Value = 1
$ws.Cells.Item(218, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(218, 4).NumberFormat = $ws.Cells.Item(219, 4).NumberFormat
